$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last existing data row (13) down into the two
# new rows (14 and 15) before writing values, so the new rows inherit the
# same number formats / styles as the rest of the table.
$ws.Range("A13:F13").Copy()
$ws.Range("A14:F15").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Row 14: 2025-11-07, 四方坪站
$ws.Range("A14").Value = 45968
$ws.Range("B14").Value = "四方坪站"
$ws.Range("C14").Value = 8017.18
$ws.Range("D14").Value = 7079.48
$ws.Range("E14").Value = 2656.04
$ws.Range("F14").Value = 383

# Row 15: 2025-11-07, 高岭站
$ws.Range("A15").Value = 45968
$ws.Range("B15").Value = "高岭站"
$ws.Range("C15").Value = 4299.03
$ws.Range("D15").Value = 3628.58
$ws.Range("E15").Value = 1119.68
$ws.Range("F15").Value = 157

# Move the active selection, matching the author's cursor position after
# the edit.
$ws.Range("J12").Select()
